$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.895.88"
$ws.Range("E2").Value = "  +0.19%  "
$ws.Range("D3").Value = "1.889.61"
$ws.Range("E3").Value = "  -0.04%  "
$ws.Range("E4").Value = "  +0.08%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.7686"
$ws.Range("E5").Value = "  -1.71%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "242.62"
$ws.Range("E6").Value = "  -0.55%  "
$ws.Range("E7").Value = "  +0.06%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3130"
$ws.Range("E8").Value = "  -0.40%  "
$ws.Range("E9").Value = "  +1.16%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.07135"
$ws.Range("E10").Value = "  -2.64%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.08535"
$ws.Range("E11").Value = "  +4.96%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.7642"
$ws.Range("E12").Value = "  -0.21%  "
$ws.Range("D13").Value = "1.901.97"
$ws.Range("E13").Value = "  +0.85%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.366"
$ws.Range("E14").Value = "  -1.84%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "93.73"
$ws.Range("E15").Value = "  +0.70%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "6.148"
$ws.Range("E16").Value = "  -0.88%  "
$ws.Range("D17").Value = "29.930.80"
$ws.Range("E17").Value = "  +0.32%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "13.78"
$ws.Range("E18").Value = "  -1.13%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "244.42"
$ws.Range("E19").Value = "  -0.40%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.000007813"
$ws.Range("E20").Value = "  -0.83%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.9993"
$ws.Range("E21").Value = "  -0.01%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "8.038"
$ws.Range("E22").Value = "  -1.26%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "1.001"
$ws.Range("E23").Value = "  +0.12%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.1629"
$ws.Range("E24").Value = "  +2.26%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "9.396"
$ws.Range("E25").Value = "  -0.65%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "163.07"
$ws.Range("E26").Value = "  +0.79%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "18.79"
$ws.Range("E27").Value = "  +0.07%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.039"
$ws.Range("E28").Value = "  +0.23%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.512"
$ws.Range("E29").Value = "  +4.14%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.535"
$ws.Range("E30").Value = "  -0.46%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.510"
$ws.Range("E31").Value = "  +0.73%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.121"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.05448"
$ws.Range("E33").Value = "  -2.63%  "
$ws.Range("E34").Value = "  -0.74%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.7456"
$ws.Range("E35").Value = "  -1.27%  "
$ws.Range("E36").Value = "  +0.52%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.696"
$ws.Range("E37").Value = "  +2.19%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.01950"
$ws.Range("E38").Value = "  +0.80%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.779"
$ws.Range("E39").Value = "  +0.04%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.4473"
$ws.Range("E40").Value = "  +0.52%  "
$ws.Range("D41").Value = "1.102.95"
$ws.Range("E41").Value = "  -3.58%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "73.19"
$ws.Range("E42").Value = "  -0.82%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "6.083"
$ws.Range("E43").Value = "  +2.04%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.8518"
$ws.Range("E44").Value = "  -0.49%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.9999"
$ws.Range("E45").Value = "  +0.04%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "103.01"
$ws.Range("E46").Value = "  +1.09%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.870"
$ws.Range("E47").Value = "  -1.58%  "
$ws.Range("E48").Value = "  +2.06%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "3.056"
$ws.Range("E49").Value = "  -1.99%  "
$ws.Range("D50").Value = "2.024.30"
$ws.Range("E50").Value = "  +1.30%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.06082"
$ws.Range("E51").Value = "  +0.32%  "
